$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N20").Value = $null
$ws.Range("H20").Value = 1990.5
$ws.Range("I20").Value = 1990.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1990.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1760.5

$ws.Range("H34").Value = 4791.375
$ws.Range("I34").Value = 1688.5
$ws.Range("K34").Value = 1688.5
$ws.Range("M34").Value = -1485.5

$ws.Range("N35").Value = $null
$ws.Range("H35").Value = 1990.5
$ws.Range("I35").Value = 1990.5
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1990.5
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -1611.5

$ws.Range("H36").Value = 4791.375
$ws.Range("I36").Value = 1688.5
$ws.Range("K36").Value = 1688.5
$ws.Range("M36").Value = -973.5

$ws.Range("H86").Value = 3034.818
$ws.Range("I86").Value = 1609.875
$ws.Range("J86").Value = 6834.6665
$ws.Range("K86").Value = 1609.875
$ws.Range("L86").Value = 6834.6665
$ws.Range("M86").Value = -486.875
$ws.Range("N86").Value = -9080.666499999999

$ws.Range("H88").Value = 18705.666
$ws.Range("I88").Value = 1817
$ws.Range("J88").Value = 27150
$ws.Range("K88").Value = 1817
$ws.Range("L88").Value = 27150
$ws.Range("M88").Value = -1411
$ws.Range("N88").Value = -27962

$ws.Range("H89").Value = 3034.818
$ws.Range("I89").Value = 1609.875
$ws.Range("J89").Value = 6834.6665
$ws.Range("K89").Value = 8049.375
$ws.Range("L89").Value = 34173.3325
$ws.Range("M89").Value = -2433.375
$ws.Range("N89").Value = -45405.3325

$ws.Range("H91").Value = 18705.666
$ws.Range("I91").Value = 1817
$ws.Range("J91").Value = 27150
$ws.Range("K91").Value = 1817
$ws.Range("L91").Value = 27150
$ws.Range("M91").Value = -413
$ws.Range("N91").Value = -29958

$ws.Range("H132").Value = 1354.25
$ws.Range("I132").Value = 1354.25
$ws.Range("K132").Value = 4062.75
$ws.Range("M132").Value = -1532.75

$ws.Range("H133").Value = 25000
$ws.Range("J133").Value = 25000
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -35120

$ws.Range("H136").Value = 25083.334
$ws.Range("J136").Value = 25083.334
$ws.Range("L136").Value = 25083.334
$ws.Range("N136").Value = -35283.334

$ws.Range("H137").Value = 1892.2433
$ws.Range("I137").Value = 1257.125
$ws.Range("K137").Value = 3771.375
$ws.Range("M137").Value = -1221.375

$ws.Range("H138").Value = 2629.859
$ws.Range("I138").Value = 1115.0476
$ws.Range("J138").Value = 4397.1387
$ws.Range("K138").Value = 3345.142800000001
$ws.Range("L138").Value = 13191.4161
$ws.Range("M138").Value = 1794.857199999999
$ws.Range("N138").Value = -23471.4161

$ws.Range("N139").Value = $null
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0

$ws.Range("H140").Value = 35000
$ws.Range("J140").Value = 35000
$ws.Range("L140").Value = 35000
$ws.Range("N140").Value = -45360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N39").Value = $null
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0

$ws.Range("H45").Value = 7428.1665
$ws.Range("J45").Value = 1881.5
$ws.Range("L45").Value = 1881.5
$ws.Range("N45").Value = -2635.5

$ws.Range("H74").Value = 1633.1025
$ws.Range("I74").Value = 1327
$ws.Range("J74").Value = 2073.125
$ws.Range("K74").Value = 1327
$ws.Range("L74").Value = 2073.125
$ws.Range("M74").Value = -453
$ws.Range("N74").Value = -3821.125

$ws.Range("H77").Value = 1633.1025
$ws.Range("I77").Value = 1327
$ws.Range("J77").Value = 2073.125
$ws.Range("K77").Value = 6635
$ws.Range("L77").Value = 10365.625
$ws.Range("M77").Value = -2267
$ws.Range("N77").Value = -19101.625

$ws.Range("H102").Value = 3369432.2
$ws.Range("I102").Value = 6175028
$ws.Range("J102").Value = 2717.8
$ws.Range("K102").Value = 6175028
$ws.Range("L102").Value = 2717.8
$ws.Range("M102").Value = -6173406
$ws.Range("N102").Value = -5961.8

$ws.Range("H110").Value = 651.3125
$ws.Range("I110").Value = 561.4
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 561.4
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 1483.6
$ws.Range("N110").Value = -6090

$ws.Range("H122").Value = 1712866.4
$ws.Range("I122").Value = 1834999.8
$ws.Range("K122").Value = 5504999.4
$ws.Range("M122").Value = -5502549.4

$ws.Range("H132").Value = 3410.0588
$ws.Range("I132").Value = 2118.0908
$ws.Range("K132").Value = 6354.2724
$ws.Range("M132").Value = -3824.2724

$ws.Range("H135").Value = 36115.332
$ws.Range("J135").Value = 36115.332
$ws.Range("L135").Value = 36115.332
$ws.Range("N135").Value = -46255.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4612

$ws.Range("H39").Value = 4050
$ws.Range("I39").Value = 1100
$ws.Range("J39").Value = 7000
$ws.Range("K39").Value = 1100
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = -709
$ws.Range("N39").Value = -7782

$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4840

$ws.Range("H44").Value = 333346240
$ws.Range("J44").Value = 333346240
$ws.Range("L44").Value = 333346240
$ws.Range("N44").Value = -333347124

$ws.Range("H49").Value = 4050
$ws.Range("I49").Value = 1100
$ws.Range("J49").Value = 7000
$ws.Range("K49").Value = 1100
$ws.Range("L49").Value = 7000
$ws.Range("M49").Value = -918
$ws.Range("N49").Value = -7364

$ws.Range("H141").Value = 78675.60000000001
$ws.Range("J141").Value = 78675.60000000001
$ws.Range("L141").Value = 78675.60000000001
$ws.Range("N141").Value = -89035.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 6250188
$ws.Range("I38").Value = 11111133
$ws.Range("J38").Value = 400.7143
$ws.Range("K38").Value = 33333399
$ws.Range("L38").Value = 1202.1429
$ws.Range("M38").Value = -33333052
$ws.Range("N38").Value = -1896.1429

$ws.Range("H140").Value = 9523.529
$ws.Range("I140").Value = 9523.529
$ws.Range("K140").Value = 28570.587
$ws.Range("M140").Value = -23390.587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5557950.5
$ws.Range("I22").Value = 22223422
$ws.Range("J22").Value = 2793.3333
$ws.Range("K22").Value = 22223422
$ws.Range("L22").Value = 2793.3333
$ws.Range("M22").Value = -22223127
$ws.Range("N22").Value = -3383.3333

$ws.Range("H27").Value = 5557950.5
$ws.Range("I27").Value = 22223422
$ws.Range("J27").Value = 2793.3333
$ws.Range("K27").Value = 22223422
$ws.Range("L27").Value = 2793.3333
$ws.Range("M27").Value = -22223315
$ws.Range("N27").Value = -3007.3333

$ws.Range("H55").Value = 16667055
$ws.Range("I55").Value = 333.84616
$ws.Range("J55").Value = 29412194
$ws.Range("K55").Value = 333.84616
$ws.Range("L55").Value = 29412194
$ws.Range("M55").Value = -160.84616
$ws.Range("N55").Value = -29412540

$ws.Range("H61").Value = 1435.1111
$ws.Range("I61").Value = 1435.1111
$ws.Range("K61").Value = 1435.1111
$ws.Range("M61").Value = -1233.1111

$ws.Range("H82").Value = 10513140
$ws.Range("I82").Value = 15699.714
$ws.Range("J82").Value = 25209556
$ws.Range("K82").Value = 15699.714
$ws.Range("L82").Value = 25209556
$ws.Range("M82").Value = -15338.714
$ws.Range("N82").Value = -25210278

$ws.Range("H85").Value = 10513140
$ws.Range("I85").Value = 15699.714
$ws.Range("J85").Value = 25209556
$ws.Range("K85").Value = 15699.714
$ws.Range("L85").Value = 25209556
$ws.Range("M85").Value = -14451.714
$ws.Range("N85").Value = -25212052

$ws.Range("H113").Value = 1435.1111
$ws.Range("I113").Value = 1435.1111
$ws.Range("K113").Value = 1435.1111
$ws.Range("M113").Value = 734.8888999999999

$ws.Range("H136").Value = 5201.1724
$ws.Range("I136").Value = 4071.568
$ws.Range("J136").Value = 8751.357
$ws.Range("K136").Value = 12214.704
$ws.Range("L136").Value = 26254.071
$ws.Range("M136").Value = -9664.704000000002
$ws.Range("N136").Value = -31354.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2733.3333
$ws.Range("J81").Value = 2733.3333
$ws.Range("L81").Value = 5466.6666
$ws.Range("N81").Value = -7588.6666

$ws.Range("H84").Value = 2733.3333
$ws.Range("J84").Value = 2733.3333
$ws.Range("L84").Value = 27333.333
$ws.Range("N84").Value = -37941.333

$ws.Range("H126").Value = 1323.5333
$ws.Range("I126").Value = 925
$ws.Range("J126").Value = 1779
$ws.Range("K126").Value = 2775
$ws.Range("L126").Value = 5337
$ws.Range("M126").Value = -305
$ws.Range("N126").Value = -10277
